$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2761.8
$ws.Range("I40").Value = 2696.3635
$ws.Range("J40").Value = 2841.7778
$ws.Range("K40").Value = 2696.3635
$ws.Range("L40").Value = 2841.7778
$ws.Range("M40").Value = -2521.3635
$ws.Range("N40").Value = -3191.7778
$ws.Range("H64").Value = 3339.2856
$ws.Range("I64").Value = 3068.182
$ws.Range("J64").Value = 4333.3335
$ws.Range("K64").Value = 3068.182
$ws.Range("L64").Value = 4333.3335
$ws.Range("M64").Value = -2820.182
$ws.Range("N64").Value = -4829.3335
$ws.Range("H67").Value = 3339.2856
$ws.Range("I67").Value = 3068.182
$ws.Range("J67").Value = 4333.3335
$ws.Range("K67").Value = 3068.182
$ws.Range("L67").Value = 4333.3335
$ws.Range("M67").Value = -2210.182
$ws.Range("N67").Value = -6049.3335
$ws.Range("H70").Value = 2287.0435
$ws.Range("I70").Value = 863.5454999999999
$ws.Range("J70").Value = 3591.9167
$ws.Range("K70").Value = 2590.6365
$ws.Range("L70").Value = 10775.7501
$ws.Range("M70").Value = -2320.6365
$ws.Range("N70").Value = -11315.7501
$ws.Range("H73").Value = 2287.0435
$ws.Range("I73").Value = 863.5454999999999
$ws.Range("J73").Value = 3591.9167
$ws.Range("K73").Value = 2590.6365
$ws.Range("L73").Value = 10775.7501
$ws.Range("M73").Value = -1654.6365
$ws.Range("N73").Value = -12647.7501
$ws.Range("H76").Value = 3638.9395
$ws.Range("I76").Value = 3003
$ws.Range("K76").Value = 3003
$ws.Range("M76").Value = -2688
$ws.Range("H79").Value = 3638.9395
$ws.Range("I79").Value = 3003
$ws.Range("K79").Value = 3003
$ws.Range("M79").Value = -1911
$ws.Range("H87").Value = 17188.72
$ws.Range("J87").Value = 17188.72
$ws.Range("L87").Value = 17188.72
$ws.Range("N87").Value = -19684.72
$ws.Range("H90").Value = 17188.72
$ws.Range("J90").Value = 17188.72
$ws.Range("L90").Value = 51566.16
$ws.Range("N90").Value = -64046.16
$ws.Range("H116").Value = 2933.8572
$ws.Range("I116").Value = 2370.3635
$ws.Range("J116").Value = 5000
$ws.Range("K116").Value = 2370.3635
$ws.Range("L116").Value = 5000
$ws.Range("M116").Value = 1071.6365
$ws.Range("N116").Value = -11884

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 4000
$ws.Range("I63").Value = 4000
$ws.Range("J63").Value = 4000
$ws.Range("K63").Value = 4000
$ws.Range("L63").Value = 4000
$ws.Range("M63").Value = -3314
$ws.Range("N63").Value = -5372
$ws.Range("H66").Value = 4000
$ws.Range("I66").Value = 4000
$ws.Range("J66").Value = 4000
$ws.Range("K66").Value = 20000
$ws.Range("L66").Value = 20000
$ws.Range("M66").Value = -16568
$ws.Range("N66").Value = -26864
$ws.Range("H104").Value = 65000
$ws.Range("J104").Value = 65000
$ws.Range("L104").Value = 65000
$ws.Range("N104").Value = -71988

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2643.3333
$ws.Range("I105").Value = 1000
$ws.Range("J105").Value = 3465
$ws.Range("K105").Value = 1000
$ws.Range("L105").Value = 3465
$ws.Range("M105").Value = 747
$ws.Range("N105").Value = -6959

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H17").Value = 1000
$ws.Range("I17").Value = 1000
$ws.Range("K17").Value = 1000
$ws.Range("M17").Value = -826
$ws.Range("H62").Value = 2300
$ws.Range("I62").Value = 2300
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 2300
$ws.Range("L62").Value = 0
$ws.Range("M62").ClearContents()
$ws.Range("N62").Value = -1676
$ws.Range("H65").Value = 2300
$ws.Range("I65").Value = 2300
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 11500
$ws.Range("L65").Value = 0
$ws.Range("M65").ClearContents()
$ws.Range("N65").Value = -8380
$ws.Range("H122").Value = 1003154
$ws.Range("I122").Value = 2618.6667
$ws.Range("J122").Value = 2503957
$ws.Range("K122").Value = 7856.000100000001
$ws.Range("L122").Value = 7511871
$ws.Range("M122").Value = -5406.000100000001
$ws.Range("N122").Value = -7516771
$ws.Range("H134").Value = 2168.3225
$ws.Range("I134").Value = 2384.1667
$ws.Range("J134").Value = 1428.2858
$ws.Range("K134").Value = 7152.500100000001
$ws.Range("L134").Value = 4284.857400000001
$ws.Range("M134").Value = -4617.500100000001
$ws.Range("N134").Value = -9354.857400000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H19").Value = 2406.923
$ws.Range("J19").Value = 2583.3333
$ws.Range("L19").Value = 7749.999899999999
$ws.Range("N19").Value = -8097.999899999999
$ws.Range("H107").Value = 486726.78
$ws.Range("J107").Value = 716991.3
$ws.Range("L107").Value = 2150973.9
$ws.Range("N107").Value = -2154813.9
$ws.Range("H132").Value = 2178.5789
$ws.Range("I132").Value = 2074.75
$ws.Range("J132").Value = 2206.2666
$ws.Range("K132").Value = 18672.75
$ws.Range("L132").Value = 19856.3994
$ws.Range("M132").Value = -16142.75
$ws.Range("N132").Value = -24916.3994
$ws.Range("H137").Value = 5971.396
$ws.Range("J137").Value = 7649.9707
$ws.Range("L137").Value = 22949.9121
$ws.Range("N137").Value = -33149.9121

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4998.673
$ws.Range("I70").Value = 4739.381
$ws.Range("J70").Value = 5174.3228
$ws.Range("K70").Value = 4739.381
$ws.Range("L70").Value = 5174.3228
$ws.Range("M70").Value = -4469.381
$ws.Range("N70").Value = -5714.3228
$ws.Range("H73").Value = 4998.673
$ws.Range("I73").Value = 4739.381
$ws.Range("J73").Value = 5174.3228
$ws.Range("K73").Value = 4739.381
$ws.Range("L73").Value = 5174.3228
$ws.Range("M73").Value = -3803.381
$ws.Range("N73").Value = -7046.3228
$ws.Range("H80").Value = 2499.0625
$ws.Range("I80").Value = 2498.9285
$ws.Range("K80").Value = 2498.9285
$ws.Range("M80").Value = -1500.9285
$ws.Range("H83").Value = 2499.0625
$ws.Range("I83").Value = 2498.9285
$ws.Range("K83").Value = 12494.6425
$ws.Range("M83").Value = -7502.6425
$ws.Range("H102").Value = 3078678.5
$ws.Range("I102").Value = 4049956.2
$ws.Range("J102").Value = 2965.3333
$ws.Range("K102").Value = 4049956.2
$ws.Range("L102").Value = 2965.3333
$ws.Range("M102").Value = -4048334.2
$ws.Range("N102").Value = -6209.3333
$ws.Range("H122").Value = 99178.16
$ws.Range("I122").Value = 141478
$ws.Range("J122").Value = 4003.5
$ws.Range("K122").Value = 424434
$ws.Range("L122").Value = 12010.5
$ws.Range("M122").Value = -421984
$ws.Range("N122").Value = -16910.5
$ws.Range("H135").Value = 40077.184
$ws.Range("J135").Value = 40077.184
$ws.Range("L135").Value = 40077.184
$ws.Range("N135").Value = -50217.184
$ws.Range("H138").Value = 45974
$ws.Range("J138").Value = 45974
$ws.Range("L138").Value = 45974
$ws.Range("N138").Value = -56254

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H30").Value = 0
$ws.Range("I30").Value = 0
$ws.Range("K30").Value = 0
$ws.Range("M30").ClearContents()
$ws.Range("H40").Value = 105800
$ws.Range("I40").Value = 173000
$ws.Range("J40").Value = 5000
$ws.Range("K40").Value = 173000
$ws.Range("L40").Value = 5000
$ws.Range("M40").Value = -172864
$ws.Range("N40").Value = -5272
